$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 00:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1007514
$ws.Range("C4").Value = 20354
$ws.Range("D4").Value = 137720
$ws.Range("E4").Value = 813170
$ws.Range("G4").Value = 1211
$ws.Range("H4").Value = 56624

# Row 8 - Alemania
$ws.Range("B8").Value = 158434
$ws.Range("C8").Value = 664
$ws.Range("E8").Value = 37873
$ws.Range("F8").Value = 2409
$ws.Range("G8").Value = 85
$ws.Range("H8").Value = 6061

# Row 45 - Chequia
$ws.Range("B45").Value = 7445
$ws.Range("C45").Value = 41
$ws.Range("E45").Value = 4396
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 223

# Row 71 - Islandia
$ws.Range("D71").Value = 1624
$ws.Range("E71").Value = 158
$ws.Range("F71").Value = 1
